# Update the LR-pairs data table with new TPM-based values.
# The target-cluster set grows from 2 (Inflammatory-Mac, Resolving-Mac)
# to 4 (ECs, Inflammatory-Mac, Neutrophils, Resolving-Mac) per sending
# cluster (FAPs, MuSCs), so the whole data block (rows 2-9) is rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Ccl21b", "Cxcr3", "ECs", 2, 0.6666666666666666, 0.126417, 0.379251, 0.733230478333749, 0.8047954435010685, 1, 0.5, 0.358189, 0.716378, 0.052978466298774, 0.03595390511601538, 0.045281178813, 0.271687072878, 0.03884542618563846, 0.02893553901343893),
    @("FAPs", "Ccl21b", "Cxcr3", "Inflammatory-Mac", 2, 0.6666666666666666, 0.126417, 0.379251, 0.733230478333749, 0.8047954435010685, 2, 0.6666666666666666, 0.574863, 1.724589, 0.08502595018806307, 0.08655445766079338, 0.07267245587100001, 0.6540521028389999, 0.062343618127175, 0.06965863314011266),
    @("FAPs", "Ccl21b", "Cxcr3", "Neutrophils", 2, 0.6666666666666666, 0.126417, 0.379251, 0.733230478333749, 0.8047954435010685, 3, 1, 5.590969, 16.772907, 0.8269404217996371, 0.8418062905306279, 0.706793528073, 6.361141752657, 0.6063379210296601, 0.677481866929586),
    @("FAPs", "Ccl21b", "Cxcr3", "Resolving-Mac", 2, 0.6666666666666666, 0.126417, 0.379251, 0.733230478333749, 0.8047954435010685, 3, 1, 0.237009, 0.7110270000000001, 0.0350551617135259, 0.03568534669256324, 0.029961966753, 0.269657700777, 0.02570351299127552, 0.02871940441793082),
    @("MuSCs", "Ccl21b", "Cxcr3", "ECs", 1, 0.5, 0.045994, 0.091988, 0.266769521666251, 0.1952045564989315, 1, 0.5, 0.358189, 0.716378, 0.052978466298774, 0.03595390511601538, 0.016474544866, 0.06589817946399999, 0.01413304011313554, 0.007018366102576448),
    @("MuSCs", "Ccl21b", "Cxcr3", "Inflammatory-Mac", 1, 0.5, 0.045994, 0.091988, 0.266769521666251, 0.1952045564989315, 2, 0.6666666666666666, 0.574863, 1.724589, 0.08502595018806307, 0.08655445766079338, 0.026440248822, 0.158641492932, 0.02268233206088807, 0.01689582452068072),
    @("MuSCs", "Ccl21b", "Cxcr3", "Neutrophils", 1, 0.5, 0.045994, 0.091988, 0.266769521666251, 0.1952045564989315, 3, 1, 5.590969, 16.772907, 0.8269404217996371, 0.8418062905306279, 0.257151028186, 1.542906169116, 0.220602500769977, 0.1643244236010419),
    @("MuSCs", "Ccl21b", "Cxcr3", "Resolving-Mac", 1, 0.5, 0.045994, 0.091988, 0.266769521666251, 0.1952045564989315, 3, 1, 0.237009, 0.7110270000000001, 0.0350551617135259, 0.03568534669256324, 0.010900991946, 0.06540595167600001, 0.00935164872225038, 0.006965942274632421)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
